# Bestelbon.xlsx cleanup:
#  - rename the original sheet "Blad1" -> "Order"
#  - add a new "Details" sheet after it, with an "Ontvanger" / "Template" header row
#  - make "Details" the active sheet/tab

$wb = $excel.ActiveWorkbook

$order = $wb.Worksheets.Item(1)
$order.Name = "Order"

# New sheet, inserted right after "Order"
$details = $wb.Worksheets.Add($null, $order)
$details.Name = "Details"

$details.Range("A1").Value = "Ontvanger"
$details.Range("B1").Value = "Template"

# Reuse the same (bold, centered) header style as the "Order" sheet's row 1
# instead of rebuilding it attribute-by-attribute (which would create new,
# redundant style entries).
$order.Range("A1").Copy()
$details.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$details.Range("B2").Select()
